$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 2872
$ws.Range("L3").Value = 2909
$ws.Range("C4").Value = 1868
$ws.Range("L4").Value = 773
$ws.Range("L6").Value = 2614
$ws.Range("C7").Value = 28412
$ws.Range("L7").Value = 9333

# Sheet 11: Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 108

# Sheet 12: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 168
$ws.Range("L7").Value = 592

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 71
$ws.Range("L7").Value = 220

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 119
$ws.Range("L7").Value = 426

# Sheet 15: West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L2").Value = 51
$ws.Range("L7").Value = 124

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 100
$ws.Range("L7").Value = 339

# Sheet 18: Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 154

# Sheet 19: Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 46

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L5").Value = 38
$ws.Range("L7").Value = 313
$ws.Range("L8").Value = 592
$ws.Range("L12").Value = 22
$ws.Range("K19").Value = 789
$ws.Range("L19").Value = 264
$ws.Range("L20").Value = 240
$ws.Range("L23").Value = 95
$ws.Range("L29").Value = 501
$ws.Range("L30").Value = 46
$ws.Range("L33").Value = 426
$ws.Range("L36").Value = 129
$ws.Range("L37").Value = 339
$ws.Range("L42").Value = 304
$ws.Range("L44").Value = 71
$ws.Range("L47").Value = 73
$ws.Range("L48").Value = 124
$ws.Range("L49").Value = 51
$ws.Range("L52").Value = 187
$ws.Range("L53").Value = 108
$ws.Range("L54").Value = 193
$ws.Range("L60").Value = 58
$ws.Range("C63").Value = 293
$ws.Range("K63").Value = 158
$ws.Range("L63").Value = 32
$ws.Range("L67").Value = 346
$ws.Range("L71").Value = 27
$ws.Range("L76").Value = 120
$ws.Range("L77").Value = 56
$ws.Range("L79").Value = 248
$ws.Range("L81").Value = 9
$ws.Range("L83").Value = 220
$ws.Range("L84").Value = 95
$ws.Range("L85").Value = 475
$ws.Range("L86").Value = 69
$ws.Range("L88").Value = 117
$ws.Range("L89").Value = 120
$ws.Range("L90").Value = 91
$ws.Range("L91").Value = 131
$ws.Range("L93").Value = 47
$ws.Range("L95").Value = 124
$ws.Range("L97").Value = 84
$ws.Range("L99").Value = 154
$ws.Range("C101").Value = 28412
$ws.Range("L101").Value = 9333

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 104
$ws.Range("L3").Value = 127
$ws.Range("L7").Value = 346

# Sheet 22: South Deering
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 95

# Sheet 23: Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 51

# Sheet 24: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L2").Value = 41
$ws.Range("L4").Value = 16
$ws.Range("L7").Value = 193

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L6").Value = 133
$ws.Range("L7").Value = 501

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L3").Value = 25
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 124

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 90
$ws.Range("L3").Value = 81
$ws.Range("K4").Value = 33
$ws.Range("L6").Value = 81
$ws.Range("K7").Value = 789
$ws.Range("L7").Value = 264

# Sheet 28: Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L2").Value = 32
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 71

# Sheet 29: River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 120

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 88
$ws.Range("L3").Value = 92
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 304

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 95

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 131

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 82
$ws.Range("L3").Value = 88
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 248

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 75
$ws.Range("L4").Value = 23
$ws.Range("L7").Value = 240

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L3").Value = 34
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 129

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 47

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 94
$ws.Range("L3").Value = 98
$ws.Range("L6").Value = 87
$ws.Range("L7").Value = 313

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L3").Value = 27
$ws.Range("L6").Value = 35

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 73

# Sheet 65: West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L2").Value = 17
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 84

# Sheet 68: United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 117

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L3").Value = 33
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 120

# Sheet 70: Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 38

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L4").Value = 38
$ws.Range("L7").Value = 69

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 91

# Sheet 78: Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 58

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 140
$ws.Range("L3").Value = 193
$ws.Range("L6").Value = 96
$ws.Range("L7").Value = 475

# Sheet 81: Oakland
$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 27

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L2").Value = 18
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 56

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L4").Value = 18
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 187

# Sheet 91: Beverly
$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 22

# Sheet 96: Sauganash,Forest Glen
$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 9
